# Updated Pride_Converter to support v2.x of the ProteomeXchange .px file format
#
# Adds a new "Sheet1" worksheet (after the existing "submission" sheet) that
# documents the key-name overrides needed for the v2.x .px format: it lists
# the old MTD key name (column B), the new key name it used to map from
# (column C) and a generated VBA "dctKeyNameOverrides.Add(...)" snippet
# (column D).

$wb = $excel.ActiveWorkbook
$submission = $wb.Worksheets.Item(1)

# --- add the new sheet right after "submission" -----------------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $submission)
$ws.Name = "Sheet1"

# --- column A: always the "MTD" key -----------------------------------------
$ws.Range("A1:A8").Value = "MTD"

# --- column B: new (v2.x) key names -----------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("B2").Value = "email"
$ws.Range("B3").Value = "affiliation"
$ws.Range("B4").Value = "title"
$ws.Range("B5").Value = "description"
$ws.Range("B6").Value = "type"
$ws.Range("B7").Value = "comment"
$ws.Range("B8").Value = "pride_login"

# --- column C: old (v1.x) key names this new key replaces -------------------
# (C7 has no old equivalent, so it is recorded as the literal text
# "-- removed --"; a leading apostrophe is needed so it is stored as text
# with the quote-prefix flag rather than being auto-formatted)
$ws.Range("C1").Value = "submitter_name"
$ws.Range("C2").Value = "submitter_email"
$ws.Range("C3").Value = "submitter_affiliation"
$ws.Range("C4").Value = "project_title"
$ws.Range("C5").Value = "project_description"
$ws.Range("C6").Value = "submission_type"
$ws.Range("C7").Formula = "'-- removed --"
$ws.Range("C8").Value = "submitter_pride_login"

# --- column D: generated dctKeyNameOverrides.Add(...) VBA line --------------
$ws.Range("D1").Formula = '="dctKeyNameOverrides.Add("""&B1&""""&","&""""&C1&""")"'
$ws.Range("D2:D8").Formula = '="dctKeyNameOverrides.Add("""&B2&""""&","&""""&C2&""")"'
$ws.Range("D1:D8").HorizontalAlignment = -4131

# --- column widths for the new sheet ----------------------------------------
$ws.Columns.Item(3).ColumnWidth = 18.1
$ws.Columns.Item(4).ColumnWidth = 25.5

# --- view/selection state ----------------------------------------------------
$ws.Range("D4").Select() | Out-Null
$submission.Activate() | Out-Null
$submission.Range("A1").Select() | Out-Null
